# Update "想去人数" (interested-count) values that changed between scrapes.
# Sheet "展览" (sheet1 / rId1) and sheet "全部类型" (sheet4 / rId4) both list
# the same events (全部类型 has one extra leading row), so the same four
# events need updating in both sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value = 9725
$ws1.Range("F10").Value = 3319
$ws1.Range("F13").Value = 25
$ws1.Range("F19").Value = 1399

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 9725
$ws4.Range("F11").Value = 3319
$ws4.Range("F14").Value = 25
$ws4.Range("F20").Value = 1399
